# Added average results to simulation
# - Adds a new "Average" summary block (LRU / SRRIP / Hawkeye rows) to the
#   bottom of the Config1 and Config2 sheets (IPC% in column H, MPKI in
#   column I), bolded like the header row; Config2's MPKI column is further
#   highlighted in bold red.
# - Config1 gets one extra trailing blank (but styled) row after the summary.
# - Re-point the active sheet/selection to the newly added data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Config1 ("Config1" = second tab / sheet2.xml) - rows 84-87
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Config1")

$c1Data = @(
    @(84, "Average", "LRU",     56.74, 14.32),
    @(85, "Average", "SRRIP",   56.85, 14.05),
    @(86, "Average", "Hawkeye", 58.52, 14.02)
)

foreach ($row in $c1Data) {
    $r = $row[0]
    $ws1.Range("A$r").Value2 = $row[1]
    $ws1.Range("B$r").Value2 = $row[2]
    $ws1.Range("H$r").Value2 = $row[3]
    $ws1.Range("I$r").Value2 = $row[4]
    $ws1.Range("A$r`:I$r").Font.Bold = $true
}

# Trailing blank (but bold-styled) row
$ws1.Range("H87:I87").Font.Bold = $true

# ---------------------------------------------------------------------
# Config2 ("Config2" = third tab / sheet3.xml) - rows 84-86
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Config2")

$c2Data = @(
    @(84, "Average", "LRU",     61.55, 17.989999999999998),
    @(85, "Average", "SRRIP",   61.66, 17.760000000000002),
    @(86, "Average", "Hawkeye", 62.89, 17.98)
)

foreach ($row in $c2Data) {
    $r = $row[0]
    $ws2.Range("A$r").Value2 = $row[1]
    $ws2.Range("B$r").Value2 = $row[2]
    $ws2.Range("H$r").Value2 = $row[3]
    $ws2.Range("A$r`:H$r").Font.Bold = $true
    $iCell = $ws2.Range("I$r")
    $iCell.Value2 = $row[4]
    $iCell.Font.Bold = $true
    $iCell.Font.Color = 255
}

# ---------------------------------------------------------------------
# View state: Config1 becomes the active tab, scrolled/selected near the
# newly added rows; Config2 is no longer the active tab but is scrolled to
# show the new rows too.
# ---------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 73
$ws1.Range("G85").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 76
$ws2.Range("H86").Select()

$ws1.Activate()
